$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sales")

function Set-TextCell {
    param($cell, [string]$text)
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
}

# Row 15's placeholder empty Reference/Actions cells are no longer needed now
# that rows are laid out for the new sortable header columns.
$ws.Cells.Item(15, 4).ClearContents()
$ws.Cells.Item(15, 8).ClearContents()

# Row 16: new sale transaction (مركز  كنز الطفوله)
Set-TextCell $ws.Cells.Item(16, 1) "مركز  كنز الطفوله"
Set-TextCell $ws.Cells.Item(16, 2) "2025-07-21"
Set-TextCell $ws.Cells.Item(16, 3) "#1: سند قبض  | Qty: 20 | Price: 20 | Total: 400 | VAT: 60"
$ws.Cells.Item(16, 5).Value = 400
$ws.Cells.Item(16, 6).Value = 60
$ws.Cells.Item(16, 7).Value = 460
$ws.Cells.Item(16, 9).Value = $false

# Row 17: new sale transaction (عجائب للاسماك)
Set-TextCell $ws.Cells.Item(17, 1) "عجائب للاسماك "
Set-TextCell $ws.Cells.Item(17, 2) "2025-07-21"
Set-TextCell $ws.Cells.Item(17, 3) "#1: فواتير بي فايف  | Qty: 10 | Price: 22 | Total: 220 | VAT: 33"
$ws.Cells.Item(17, 5).Value = 220
$ws.Cells.Item(17, 6).Value = 33
$ws.Cells.Item(17, 7).Value = 253
$ws.Cells.Item(17, 9).Value = $false
